# Commit: "added figure and data"
# - Inserts a new "CORSIA" worksheet between "WayPoint2050" and "Swiss",
#   populated with plotted (x,y) series for net/tech/corsia/baseline curves.
# - Moves the active selection on "WayPoint2050" away from its old A2:B18 range.
# - The newly inserted sheet becomes the active tab with cell L16 selected,
#   matching Excel's default behaviour right after adding + populating a sheet.

$wb = $excel.ActiveWorkbook

# --- WayPoint2050: move selection (was A2:B18, tabSelected) to plain F32 ---
$wayPoint = $wb.Worksheets.Item("WayPoint2050")
$wayPoint.Range("F32").Select() | Out-Null

# --- Insert a new "CORSIA" sheet right after WayPoint2050, before Swiss ---
$ws = $wb.Worksheets.Add($null, $wayPoint)
$ws.Name = "CORSIA"

$ws.Cells.Item(1,1).Value = "net (x)"
$ws.Cells.Item(1,2).Value = "net (y)"
$ws.Cells.Item(1,3).Value = "tech (x)"
$ws.Cells.Item(1,4).Value = "tech (y)"
$ws.Cells.Item(1,5).Value = "corsia (x)"
$ws.Cells.Item(1,6).Value = "corsia (y)"
$ws.Cells.Item(1,7).Value = "baseline (x)"
$ws.Cells.Item(1,8).Value = "baseline (y)"
$ws.Cells.Item(2,1).Value = 2010.2642493997801
$ws.Cells.Item(2,2).Value = 456.59680200789802
$ws.Cells.Item(2,3).Value = 2012.9258356734199
$ws.Cells.Item(2,4).Value = 516.42488843541605
$ws.Cells.Item(2,5).Value = 2010.8551686656299
$ws.Cells.Item(2,6).Value = 469.349926605429
$ws.Cells.Item(2,7).Value = 2010.8551686656299
$ws.Cells.Item(2,8).Value = 469.349926605429
$ws.Cells.Item(3,1).Value = 2011.4103379266801
$ws.Cells.Item(3,2).Value = 485.56239764919502
$ws.Cells.Item(3,3).Value = 2014.0204010801699
$ws.Cells.Item(3,4).Value = 544.44005462721896
$ws.Cells.Item(3,5).Value = 2011.9529900550799
$ws.Cells.Item(3,6).Value = 493.16137823063201
$ws.Cells.Item(3,7).Value = 2011.9529900550799
$ws.Cells.Item(3,8).Value = 493.16137823063201
$ws.Cells.Item(4,1).Value = 2012.50164379783
$ws.Cells.Item(4,2).Value = 515.901056283755
$ws.Cells.Item(4,3).Value = 2015.1180976533699
$ws.Cells.Item(4,4).Value = 573.51417570658805
$ws.Cells.Item(4,5).Value = 2013.05087924596
$ws.Cells.Item(4,6).Value = 514.114095831349
$ws.Cells.Item(4,7).Value = 2013.05087924596
$ws.Cells.Item(4,8).Value = 514.114095831349
$ws.Cells.Item(5,1).Value = 2013.59929876561
$ws.Cells.Item(5,2).Value = 546.72940051451303
$ws.Cells.Item(5,3).Value = 2016.17918046721
$ws.Cells.Item(5,4).Value = 602.62574936779902
$ws.Cells.Item(5,5).Value = 2014.14871373342
$ws.Cells.Item(5,6).Value = 537.37329202000296
$ws.Cells.Item(5,7).Value = 2014.14871373342
$ws.Cells.Item(5,8).Value = 537.37329202000296
$ws.Cells.Item(6,1).Value = 2014.6969337011501
$ws.Cells.Item(6,2).Value = 578.40237070705098
$ws.Cells.Item(6,3).Value = 2017.2402344434299
$ws.Cells.Item(6,4).Value = 632.953213149813
$ws.Cells.Item(6,5).Value = 2015.24654282759
$ws.Cells.Item(6,6).Value = 560.85988750605998
$ws.Cells.Item(6,7).Value = 2015.24654282759
$ws.Cells.Item(6,8).Value = 560.85988750605998
$ws.Cells.Item(7,1).Value = 2015.79455476822
$ws.Cells.Item(7,2).Value = 610.66008195005202
$ws.Cells.Item(7,3).Value = 2018.33785474003
$ws.Cells.Item(7,4).Value = 665.24341000672996
$ws.Cells.Item(7,5).Value = 2016.34438424928
$ws.Cells.Item(7,6).Value = 583.82671316948199
$ws.Cells.Item(7,7).Value = 2016.34438424928
$ws.Cells.Item(7,8).Value = 583.82671316948199
$ws.Cells.Item(8,1).Value = 2016.8921627372799
$ws.Cells.Item(8,2).Value = 643.47004862960102
$ws.Cells.Item(8,3).Value = 2019.4354511520301
$ws.Cells.Item(8,4).Value = 698.54066089499804
$ws.Cells.Item(8,5).Value = 2017.44222567098
$ws.Cells.Item(8,6).Value = 606.79353883290503
$ws.Cells.Item(8,7).Value = 2017.44222567098
$ws.Cells.Item(8,8).Value = 606.79353883290503
$ws.Cells.Item(9,1).Value = 2017.9897560674101
$ws.Cells.Item(9,2).Value = 676.89724197352803
$ws.Cells.Item(9,3).Value = 2020.53302907274
$ws.Cells.Item(9,4).Value = 732.61756651721703
$ws.Cells.Item(9,5).Value = 2018.5400678631399
$ws.Cells.Item(9,6).Value = 629.72787888241396
$ws.Cells.Item(9,7).Value = 2018.5400678631399
$ws.Cells.Item(9,8).Value = 629.72787888241396
$ws.Cells.Item(10,1).Value = 2019.08732705388
$ws.Cells.Item(10,2).Value = 711.26651812097896
$ws.Cells.Item(10,3).Value = 2021.6122915759699
$ws.Cells.Item(10,4).Value = 767.07318063972104
$ws.Cells.Item(10,5).Value = 2019.6379285466101
$ws.Cells.Item(10,6).Value = 651.88256419797096
$ws.Cells.Item(10,7).Value = 2019.6379285466101
$ws.Cells.Item(10,8).Value = 651.88256419797096
$ws.Cells.Item(11,1).Value = 2020.18488417188
$ws.Cells.Item(11,2).Value = 746.22053531889196
$ws.Cells.Item(11,3).Value = 2022.6915343923399
$ws.Cells.Item(11,4).Value = 802.35885820742203
$ws.Cells.Item(11,5).Value = 2019.8886628077501
$ws.Cells.Item(11,6).Value = 659.24614365530294
$ws.Cells.Item(11,7).Value = 2019.8886628077501
$ws.Cells.Item(11,8).Value = 659.24614365530294
$ws.Cells.Item(12,1).Value = 2021.28241663481
$ws.Cells.Item(12,2).Value = 782.21409216207303
$ws.Cells.Item(12,3).Value = 2023.7890483639801
$ws.Cells.Item(12,4).Value = 839.13206978455298
$ws.Cells.Item(12,5).Value = 2020.83451886178
$ws.Cells.Item(12,6).Value = 678.80205261789797
$ws.Cells.Item(12,7).Value = 2021.49476760121
$ws.Cells.Item(12,8).Value = 659.24614365530294
$ws.Cells.Item(13,1).Value = 2022.37992521315
$ws.Cells.Item(13,2).Value = 819.21470303660703
$ws.Cells.Item(13,3).Value = 2024.8865407624301
$ws.Cells.Item(13,4).Value = 876.81487855129501
$ws.Cells.Item(13,5).Value = 2021.52245601807
$ws.Cells.Item(13,6).Value = 695.53803375829705
$ws.Cells.Item(13,7).Value = 2022.5464060075201
$ws.Cells.Item(13,8).Value = 659.24614365530294
$ws.Cells.Item(14,1).Value = 2023.4591170721901
$ws.Cells.Item(14,2).Value = 856.64891189700597
$ws.Cells.Item(14,3).Value = 2025.9474328297399
$ws.Cells.Item(14,4).Value = 913.96896205736505
$ws.Cells.Item(14,5).Value = 2022.62022039269
$ws.Cells.Item(14,6).Value = 721.75342081318104
$ws.Cells.Item(14,7).Value = 2023.5980444138199
$ws.Cells.Item(14,8).Value = 659.24614365530294
$ws.Cells.Item(15,1).Value = 2024.50169322291
$ws.Cells.Item(15,2).Value = 894.20274753996102
$ws.Cells.Item(15,3).Value = 2026.9077059491101
$ws.Cells.Item(15,4).Value = 948.321395347921
$ws.Cells.Item(15,5).Value = 2023.7179716693099
$ws.Cells.Item(15,6).Value = 748.52106330461299
$ws.Cells.Item(15,7).Value = 2024.64968282012
$ws.Cells.Item(15,8).Value = 659.24614365530294
$ws.Cells.Item(16,1).Value = 2025.52594351874
$ws.Cells.Item(16,2).Value = 932.57536773376899
$ws.Cells.Item(16,3).Value = 2027.8222438754899
$ws.Cells.Item(16,4).Value = 981.37115569531397
$ws.Cells.Item(16,5).Value = 2024.8157167821701
$ws.Cells.Item(16,6).Value = 775.548590707362
$ws.Cells.Item(16,7).Value = 2025.7013212264301
$ws.Cells.Item(16,8).Value = 659.24614365530294
$ws.Cells.Item(17,1).Value = 2026.5135706952599
$ws.Cells.Item(17,2).Value = 971.38008717076605
$ws.Cells.Item(17,3).Value = 2028.7459068934299
$ws.Cells.Item(17,4).Value = 1015.60681094837
$ws.Cells.Item(17,5).Value = 2025.9134457151399
$ws.Cells.Item(17,6).Value = 803.258316002318
$ws.Cells.Item(17,7).Value = 2026.7529596327299
$ws.Cells.Item(17,8).Value = 659.24614365530294
$ws.Cells.Item(18,1).Value = 2027.4279898709499
$ws.Cells.Item(18,2).Value = 1009.4367726343399
$ws.Cells.Item(18,3).Value = 2029.6878324899201
$ws.Cells.Item(18,4).Value = 1051.6916235650999
$ws.Cells.Item(18,5).Value = 2027.0111654024599
$ws.Cells.Item(18,6).Value = 831.35786866424803
$ws.Cells.Item(18,7).Value = 2027.80459803903
$ws.Cells.Item(18,8).Value = 659.24614365530294
$ws.Cells.Item(19,1).Value = 2028.2234929654201
$ws.Cells.Item(19,2).Value = 1044.2995937071901
$ws.Cells.Item(19,3).Value = 2030.51332106167
$ws.Cells.Item(19,4).Value = 1084.3627385935999
$ws.Cells.Item(19,5).Value = 2028.10886428707
$ws.Cells.Item(19,6).Value = 860.33453290187401
$ws.Cells.Item(19,7).Value = 2028.8562364453401
$ws.Cells.Item(19,8).Value = 659.24614365530294
$ws.Cells.Item(20,1).Value = 2028.99153943186
$ws.Cells.Item(20,2).Value = 1079.03485885745
$ws.Cells.Item(20,3).Value = 2031.32497632645
$ws.Cells.Item(20,4).Value = 1118.2004711719201
$ws.Cells.Item(20,5).Value = 2029.2065454508599
$ws.Cells.Item(20,6).Value = 890.05836625953498
$ws.Cells.Item(20,7).Value = 2029.9078748516399
$ws.Cells.Item(20,8).Value = 659.24614365530294
$ws.Cells.Item(21,1).Value = 2029.8052852539299
$ws.Cells.Item(21,2).Value = 1116.5838191154401
$ws.Cells.Item(21,3).Value = 2031.9020046782
$ws.Cells.Item(21,4).Value = 1149.28707951507
$ws.Cells.Item(21,5).Value = 2030.30419425487
$ws.Cells.Item(21,6).Value = 921.14659540161006
$ws.Cells.Item(21,7).Value = 2030.95951325794
$ws.Cells.Item(21,8).Value = 659.24614365530294
$ws.Cells.Item(22,1).Value = 2030.5458757403001
$ws.Cells.Item(22,2).Value = 1151.16420782515
$ws.Cells.Item(22,3).Value = 2032.6871155853401
$ws.Cells.Item(22,4).Value = 1177.0151144517899
$ws.Cells.Item(22,5).Value = 2031.4017860440499
$ws.Cells.Item(22,6).Value = 954.63875997336595
$ws.Cells.Item(22,7).Value = 2032.0111516642501
$ws.Cells.Item(22,8).Value = 659.24614365530294
$ws.Cells.Item(23,1).Value = 2031.2772820963501
$ws.Cells.Item(23,2).Value = 1187.04797277433
$ws.Cells.Item(23,3).Value = 2033.4643031491701
$ws.Cells.Item(23,4).Value = 1212.2614151535599
$ws.Cells.Item(23,5).Value = 2032.4993084908699
$ws.Cells.Item(23,6).Value = 991.05462979743697
$ws.Cells.Item(23,7).Value = 2033.0627900705499
$ws.Cells.Item(23,8).Value = 659.24614365530294
$ws.Cells.Item(24,1).Value = 2031.99953777679
$ws.Cells.Item(24,2).Value = 1222.8245516433999
$ws.Cells.Item(24,3).Value = 2034.2506036433799
$ws.Cells.Item(24,4).Value = 1249.2063665226401
$ws.Cells.Item(24,5).Value = 2033.5419122982
$ws.Cells.Item(24,6).Value = 1027.4423705914901
$ws.Cells.Item(24,7).Value = 2034.11442847685
$ws.Cells.Item(24,8).Value = 659.24614365530294
$ws.Cells.Item(25,1).Value = 2032.56635372711
$ws.Cells.Item(25,2).Value = 1251.6453616388701
$ws.Cells.Item(25,3).Value = 2035.0368801229199
$ws.Cells.Item(25,4).Value = 1287.1638565072401
$ws.Cells.Item(25,5).Value = 2034.5204624017299
$ws.Cells.Item(25,6).Value = 1063.0365723653499
$ws.Cells.Item(25,7).Value = 2035.1660668831601
$ws.Cells.Item(25,8).Value = 659.24614365530294
$ws.Cells.Item(26,1).Value = 2033.14230751893
$ws.Cells.Item(26,2).Value = 1281.1144861908299
$ws.Cells.Item(26,3).Value = 2035.7865603538
$ws.Cells.Item(26,4).Value = 1324.4204896665301
$ws.Cells.Item(26,5).Value = 2035.3983977590401
$ws.Cells.Item(26,6).Value = 1095.6519655089701
$ws.Cells.Item(26,7).Value = 2036.2177052894599
$ws.Cells.Item(26,8).Value = 659.24614365530294
$ws.Cells.Item(27,1).Value = 2033.8370774590401
$ws.Cells.Item(27,2).Value = 1317.9909808260199
$ws.Cells.Item(27,3).Value = 2036.4722179897501
$ws.Cells.Item(27,4).Value = 1359.57192971425
$ws.Cells.Item(27,5).Value = 2036.2671575285799
$ws.Cells.Item(27,6).Value = 1129.2105524812
$ws.Cells.Item(27,7).Value = 2037.26934369576
$ws.Cells.Item(27,8).Value = 659.24614365530294
$ws.Cells.Item(28,1).Value = 2034.45195595532
$ws.Cells.Item(28,2).Value = 1351.5460756892101
$ws.Cells.Item(28,3).Value = 2037.13956383645
$ws.Cells.Item(28,4).Value = 1394.9491002468301
$ws.Cells.Item(28,5).Value = 2037.11761010895
$ws.Cells.Item(28,6).Value = 1162.80091492468
$ws.Cells.Item(28,7).Value = 2038.32098210207
$ws.Cells.Item(28,8).Value = 659.24614365530294
$ws.Cells.Item(29,1).Value = 2035.11684287657
$ws.Cells.Item(29,2).Value = 1389.03916807059
$ws.Cells.Item(29,3).Value = 2037.8343147242899
$ws.Cells.Item(29,4).Value = 1432.62890201502
$ws.Cells.Item(29,5).Value = 2037.9314589678099
$ws.Cells.Item(29,6).Value = 1196.00550056308
$ws.Cells.Item(29,7).Value = 2039.3726205083699
$ws.Cells.Item(29,8).Value = 659.24614365530294
$ws.Cells.Item(30,1).Value = 2036.30511005433
$ws.Cells.Item(30,2).Value = 1458.65641841044
$ws.Cells.Item(30,3).Value = 2038.5107511851199
$ws.Cells.Item(30,4).Value = 1470.6456506806201
$ws.Cells.Item(30,5).Value = 2038.8092938228999
$ws.Cells.Item(30,6).Value = 1232.85840256914
$ws.Cells.Item(30,7).Value = 2040.42425891467
$ws.Cells.Item(30,8).Value = 659.24614365530294
$ws.Cells.Item(31,1).Value = 2036.8352198555599
$ws.Cells.Item(31,2).Value = 1491.4101398867001
$ws.Cells.Item(31,3).Value = 2039.15059275148
$ws.Cells.Item(31,4).Value = 1507.9044452872599
$ws.Cells.Item(31,5).Value = 2039.62308731592
$ws.Cells.Item(31,6).Value = 1268.3973965001301
$ws.Cells.Item(32,1).Value = 2037.2647733654201
$ws.Cells.Item(32,2).Value = 1518.72040146731
$ws.Cells.Item(32,3).Value = 2039.7172927443701
$ws.Cells.Item(32,4).Value = 1541.6144066740401
$ws.Cells.Item(33,1).Value = 2037.6943088882899
$ws.Cells.Item(33,2).Value = 1546.7890544709401
$ws.Cells.Item(34,1).Value = 2038.1969233070799
$ws.Cells.Item(34,2).Value = 1581.0492350459101
$ws.Cells.Item(35,1).Value = 2038.72691867039
$ws.Cells.Item(35,2).Value = 1618.62804053799
$ws.Cells.Item(36,1).Value = 2039.2203101719299
$ws.Cells.Item(36,2).Value = 1655.82698341779
$ws.Cells.Item(37,1).Value = 2039.6770923337299
$ws.Cells.Item(37,2).Value = 1692.8770328303799

$ws.Range("L16").Select() | Out-Null
